# Region IX_GABALDON.xlsx edit:
#  - Insert a new first column "INDEX (DO NOT MODIFY)" (numbered index per row)
#  - Rename several header labels (uppercase some + tweak "Region" -> "REGION")
#  - dataValidation sqref + dimension auto-update because of the insert

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new column before column A, shifting everything right.
$ws.Columns.Item(1).Insert("ShiftToRight")

# 2. Give the new column A the same look (styles) as column B (the former column A).
$ws.Range("B1:B10").Copy()
$ws.Range("A1:A10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match the width used in the diff for the brand-new column (23 chars).
$ws.Columns.Item(1).ColumnWidth = 22.17

# 3. Fill in the new INDEX column.
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$ws.Range("A2").Value = 50
$ws.Range("A3").Value = 51
$ws.Range("A4").Value = 147
$ws.Range("A5").Value = 194
$ws.Range("A6").Value = 195
$ws.Range("A7").Value = 246
$ws.Range("A8").Value = 247
$ws.Range("A9").Value = 248
$ws.Range("A10").Value = 249

# 4. Update header text (row 1) that changed wording/casing.
$ws.Range("C1").Value = "REGION"
$ws.Range("R1").Value = " TARGET COMPLETION DATE "
$ws.Range("S1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("T1").Value = "PROJECT ID"
$ws.Range("U1").Value = "CONTRACT ID"
$ws.Range("V1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("W1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("X1").Value = "BID OPENING"
$ws.Range("Y1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("Z1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("AA1").Value = "NAME OF CONTRACTOR"
$ws.Range("AB1").Value = "OTHER REMARKS"
